$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45, shifting existing rows 45.. down by one.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new data record.
$ws.Range("A45").Value = 2
$ws.Range("B45").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C45").Value = "Coquimbo"
$ws.Range("D45").Value = 45203
$ws.Range("E45").Value = 4
$ws.Range("F45").Value = 100112026
$ws.Range("G45").Value = "Haba"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 1300
$ws.Range("K45").Value = 7000
$ws.Range("L45").Value = 8000
$ws.Range("M45").Value = 7500
$ws.Range("N45").Value = "$/saco 25 kilos"
$ws.Range("O45").Value = "Provincia de Limarí"
$ws.Range("P45").Value = 300
$ws.Range("Q45").Value = 25
$ws.Range("R45").Value = "Hortaliza"
